# Add presenter notes to slide 2 and slide 3.
#
# PowerPoint stores speaker notes as a "notes page" attached to each
# slide (Slide.NotesPage). The notes page carries a body placeholder
# (ppPlaceholderBody / type 2) whose TextFrame holds the actual text -
# exactly like the main slide body, just on the notes page instead.

$p = $ppt.ActivePresentation

# --- Slide 2: "Some notes on the second slide." ---------------------
$slide2 = $p.Slides.Item(2)
$notes2 = $slide2.NotesPage
$notesBody2 = $notes2.Shapes.AddPlaceholder(2)
$notesBody2.TextFrame.TextRange.Text = "Some notes on the second slide."

# --- Slide 3: two lines of notes -------------------------------------
$slide3 = $p.Slides.Item(3)
$notes3 = $slide3.NotesPage
$notesBody3 = $notes3.Shapes.AddPlaceholder(2)
$notesBody3.TextFrame.TextRange.Text = "Final notes on the third slide.`nSecond line of notes."
